$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 7.107000000000001
$ws.Range("C5").Value = -13.155
$ws.Range("A8").Value = -21.532
$ws.Range("C8").Value = -12.891
$ws.Range("A10").Value = -21.464
$ws.Range("B11").Value = 7.142
$ws.Range("A12").Value = -21.303
$ws.Range("B12").Value = 6.304
$ws.Range("C12").Value = -11.214
$ws.Range("C13").Value = -12.981
$ws.Range("B15").Value = 5.01
$ws.Range("C15").Value = -13.132
$ws.Range("B17").Value = 4.809
$ws.Range("A18").Value = -21.839
$ws.Range("C21").Value = -13.122
$ws.Range("A25").Value = -21.779
$ws.Range("C25").Value = -12.181
$ws.Range("B26").Value = 6.476999999999999
$ws.Range("B27").Value = 5.520999999999999
$ws.Range("B28").Value = 4.737
$ws.Range("B32").Value = 6.642999999999999
$ws.Range("C32").Value = -12.564
$ws.Range("C36").Value = -13.173
$ws.Range("A37").Value = -21.3
$ws.Range("B37").Value = 6.101000000000001
$ws.Range("C38").Value = -11.859
$ws.Range("B41").Value = 7.340999999999999
$ws.Range("C41").Value = -12.758
$ws.Range("B47").Value = 4.846
$ws.Range("C50").Value = -13.026
$ws.Range("B51").Value = 6.132
$ws.Range("C52").Value = -11.951
$ws.Range("A55").Value = -21.767
$ws.Range("C59").Value = -12.266
$ws.Range("B65").Value = 6.348000000000001
$ws.Range("C67").Value = -11.065
$ws.Range("A68").Value = -21.682
$ws.Range("B73").Value = 6.792
$ws.Range("A77").Value = -21.042
$ws.Range("A78").Value = -20.909
$ws.Range("A79").Value = -21.128
$ws.Range("A80").Value = -20.887
$ws.Range("A81").Value = -21.501
$ws.Range("A82").Value = -21.843
$ws.Range("A84").Value = -21.67
$ws.Range("B84").Value = 6.203
$ws.Range("C84").Value = -12.344
$ws.Range("B85").Value = 5.126
$ws.Range("C86").Value = -13.31
$ws.Range("C88").Value = -13.313
$ws.Range("B89").Value = 4.55
$ws.Range("C89").Value = -13.317
$ws.Range("B93").Value = 6.077000000000001
$ws.Range("B95").Value = 7.145999999999999
$ws.Range("C95").Value = -11.799
$ws.Range("B98").Value = 7.324000000000001
$ws.Range("B99").Value = 5.894
$ws.Range("A101").Value = -21.851
$ws.Range("B101").Value = 5.427000000000001
$ws.Range("A102").Value = -20.576
$ws.Range("B102").Value = 7.222000000000001
$ws.Range("C105").Value = -12.653
